# Updated cryptos list with latest price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.891.56'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '3.035.12'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'594.03"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").Value = "'153.74"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +5.93%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '3.033.21'
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("D10").Value = "'6.78"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +12.26%  '
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").Value = "'35.60"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +3.08%  '
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").Value = '3.537.63'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '62.854.24'
$ws.Range("E18").Value = '  +2.01%  '
$ws.Range("D19").Value = '3.035.02'
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = "'453.45"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = "'14.27"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").Value = "'7.51"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("D24").Value = "'83.16"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E25").Value = '  +2.72%  '
$ws.Range("D26").Value = "'2.30"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("D27").Value = "'12.29"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +2.06%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = "'7.53"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +3.30%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = "'2.70"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'2.24"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +7.89%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = "'27.56"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").Value = '0.0₃0867'
$ws.Range("E35").Value = '  +3.81%  '
$ws.Range("E36").Value = '  +1.85%  '
$ws.Range("D38").Value = "'3.20"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +10.56%  '
$ws.Range("D39").Value = "'2.12"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +1.69%  '
$ws.Range("E40").Value = '  +4.75%  '
$ws.Range("D41").Value = "'50.39"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").Value = "'9.09"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("D43").Value = "'0.303"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +10.65%  '
$ws.Range("D44").Value = "'42.64"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +7.25%  '
$ws.Range("D45").Value = "'394.49"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").Value = "'0.0361"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +1.76%  '
$ws.Range("D47").Value = '2.725.81'
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = "'132.26"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("D49").Value = "'2.32"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +7.37%  '
$ws.Range("D51").Value = "'24.45"
$ws.Range("D51").Style = 'Normal'
